# Fruta / hortaliza, semanal
#
# Two new weekly price records (Kiwi, Terminal Hortofrutícola Agro Chillán)
# are inserted as new rows 154 and 155. The previously-existing rows 154..256
# shift down, unchanged, to become rows 156..258.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows at position 154; everything currently at/after row 154
# (through row 256) moves down to make room, ending at row 258.
$ws.Range("A154:A155").EntireRow.Insert()

# ---- New row 154 ----
$row = New-Object 'object[,]' 1,20
$row[0,0]  = 7
$row[0,1]  = "Terminal Hortofrutícola Agro Chillán"
$row[0,2]  = "Ñuble"
$row[0,3]  = 45068
$row[0,4]  = 16
$row[0,5]  = "Fruta"
$row[0,6]  = 100101
$row[0,7]  = "Berries"
$row[0,8]  = 100101007
$row[0,9]  = "Kiwi"
$row[0,10] = "Hayward"
$row[0,11] = "Especial"
$row[0,12] = 50
$row[0,13] = 12000
$row[0,14] = 12000
$row[0,15] = 12000
$row[0,16] = "`$/bandeja 18 kilos"
$row[0,17] = "Región de O'Higgins"
$row[0,18] = 667
$row[0,19] = 18
$ws.Range("A154:T154").Value = $row

# ---- New row 155 ----
$row = New-Object 'object[,]' 1,20
$row[0,0]  = 7
$row[0,1]  = "Terminal Hortofrutícola Agro Chillán"
$row[0,2]  = "Ñuble"
$row[0,3]  = 45068
$row[0,4]  = 16
$row[0,5]  = "Fruta"
$row[0,6]  = 100101
$row[0,7]  = "Berries"
$row[0,8]  = 100101007
$row[0,9]  = "Kiwi"
$row[0,10] = "Hayward"
$row[0,11] = "Primera"
$row[0,12] = 80
$row[0,13] = 10000
$row[0,14] = 10000
$row[0,15] = 10000
$row[0,16] = "`$/bandeja 18 kilos"
$row[0,17] = "Región de O'Higgins"
$row[0,18] = 556
$row[0,19] = 18
$ws.Range("A155:T155").Value = $row

# D column keeps its date style (s="2") from the row above/below after an
# insert, but set it explicitly to be safe.
$ws.Range("D154:D155").NumberFormat = $ws.Range("D153").NumberFormat
